$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "70.879.50"
$ws.Range("E2").Value = "  +4.34%  "

$ws.Range("D3").Value = "3.558.39"
$ws.Range("E3").Value = "  +3.44%  "

$ws.Range("E4").Value = "  +0.13%  "

$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = "597.75"
$c.Style = "Normal"
$ws.Range("E5").Value = "  +3.27%  "

$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = "171.83"
$c.Style = "Normal"
$ws.Range("E6").Value = "  +3.98%  "

$ws.Range("D7").Value = "3.549.81"
$ws.Range("E7").Value = "  +3.49%  "

$ws.Range("E8").Value = "  +2.20%  "

$ws.Range("E9").Value = "  +0.03%  "

$c = $ws.Range("D10")
$c.NumberFormat = "@"
$c.Value = "0.195"
$c.Style = "Normal"
$ws.Range("E10").Value = "  +6.37%  "

$c = $ws.Range("D11")
$c.NumberFormat = "@"
$c.Value = "7.40"
$c.Style = "Normal"
$ws.Range("E11").Value = "  +9.99%  "

$c = $ws.Range("D12")
$c.NumberFormat = "@"
$c.Value = "0.586"
$c.Style = "Normal"
$ws.Range("E12").Value = "  +3.17%  "

$c = $ws.Range("D13")
$c.NumberFormat = "@"
$c.Value = "46.27"
$c.Style = "Normal"
$ws.Range("E13").Value = "  +0.24%  "

$c = $ws.Range("D14")
$c.NumberFormat = "@"
$c.Value = "0.0000276"
$c.Style = "Normal"
$ws.Range("E14").Value = "  +2.53%  "

$ws.Range("D15").Value = "4.133.94"
$ws.Range("E15").Value = "  +3.23%  "

$ws.Range("E16").Value = "  +0.76%  "

$c = $ws.Range("D17")
$c.NumberFormat = "@"
$c.Value = "611.48"
$c.Style = "Normal"
$ws.Range("E17").Value = "  +0.39%  "

$ws.Range("B18").Value = "WrappedEther"
$ws.Range("C18").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D18").Value = "3.565.43"
$ws.Range("E18").Value = "  +3.21%  "

$ws.Range("B19").Value = "WrappedBTC"
$ws.Range("C19").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D19").Value = "70.853.69"
$ws.Range("E19").Value = "  +4.26%  "

$ws.Range("E20").Value = "  +0.93%  "

$c = $ws.Range("D21")
$c.NumberFormat = "@"
$c.Value = "17.33"
$c.Style = "Normal"
$ws.Range("E21").Value = "  +1.27%  "

$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value = "0.878"
$c.Style = "Normal"
$ws.Range("E22").Value = "  +0.99%  "

$c = $ws.Range("D23")
$c.NumberFormat = "@"
$c.Value = "9.23"
$c.Style = "Normal"
$ws.Range("E23").Value = "  -15.01%  "

$c = $ws.Range("D24")
$c.NumberFormat = "@"
$c.Value = "15.67"
$c.Style = "Normal"
$ws.Range("E24").Value = "  +1.95%  "

$ws.Range("E25").Value = "  +1.74%  "

$c = $ws.Range("D26")
$c.NumberFormat = "@"
$c.Value = "3.70"
$c.Style = "Normal"
$ws.Range("E26").Value = "  -1.06%  "

$c = $ws.Range("D27")
$c.NumberFormat = "@"
$c.Value = "0.999"
$c.Style = "Normal"
$ws.Range("E27").Value = "  -0.25%  "

$c = $ws.Range("D28")
$c.NumberFormat = "@"
$c.Value = "2.60"
$c.Style = "Normal"
$ws.Range("E28").Value = "  +1.52%  "

$c = $ws.Range("D29")
$c.NumberFormat = "@"
$c.Value = "33.98"
$c.Style = "Normal"
$ws.Range("E29").Value = "  +6.01%  "

$c = $ws.Range("D30")
$c.NumberFormat = "@"
$c.Value = "9.03"
$c.Style = "Normal"
$ws.Range("E30").Value = "  +0.98%  "

$c = $ws.Range("D31")
$c.NumberFormat = "@"
$c.Value = "688.91"
$c.Style = "Normal"
$ws.Range("E31").Value = "  +13.57%  "

$c = $ws.Range("D32")
$c.NumberFormat = "@"
$c.Value = "3.06"
$c.Style = "Normal"
$ws.Range("E32").Value = "  +0.68%  "

$c = $ws.Range("D33")
$c.NumberFormat = "@"
$c.Value = "8.23"
$c.Style = "Normal"
$ws.Range("E33").Value = "  -1.07%  "

$c = $ws.Range("D34")
$c.NumberFormat = "@"
$c.Value = "7.06"
$c.Style = "Normal"
$ws.Range("E34").Value = "  +4.21%  "

$ws.Range("E35").Value = "  +0.66%  "

$ws.Range("E36").Value = "  +5.79%  "

$c = $ws.Range("D37")
$c.NumberFormat = "@"
$c.Value = "0.100"
$c.Style = "Normal"
$ws.Range("E37").Value = "  +0.44%  "

$ws.Range("E38").Value = "  +1.84%  "

$c = $ws.Range("D39")
$c.NumberFormat = "@"
$c.Value = "0.0477"
$c.Style = "Normal"
$ws.Range("E39").Value = "  +10.78%  "

$c = $ws.Range("D40")
$c.NumberFormat = "@"
$c.Value = "56.90"
$c.Style = "Normal"
$ws.Range("E40").Value = "  +1.01%  "

$c = $ws.Range("D41")
$c.NumberFormat = "@"
$c.Value = "1.00"
$c.Style = "Normal"
$ws.Range("E41").Value = "  +0.09%  "

$c = $ws.Range("D42")
$c.NumberFormat = "@"
$c.Value = "0.143"
$c.Style = "Normal"
$ws.Range("E42").Value = "  +7.00%  "

$ws.Range("D43").Value = "3.366.38"
$ws.Range("E43").Value = "  +0.38%  "

$c = $ws.Range("D44")
$c.NumberFormat = "@"
$c.Value = "0.317"
$c.Style = "Normal"
$ws.Range("E44").Value = "  -0.62%  "

$ws.Range("D45").Value = "0.0₃0698"
$ws.Range("E45").Value = "  +2.33%  "

$ws.Range("B46").Value = "InjectiveProtocol"
$ws.Range("C46").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$c = $ws.Range("D46")
$c.NumberFormat = "@"
$c.Value = "32.48"
$c.Style = "Normal"
$ws.Range("E46").Value = "  +1.19%  "

$ws.Range("B47").Value = "ThetaToken"
$ws.Range("C47").Value = "https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta"
$c = $ws.Range("D47")
$c.NumberFormat = "@"
$c.Value = "2.92"
$c.Style = "Normal"
$ws.Range("E47").Value = "  +8.74%  "

$c = $ws.Range("D48")
$c.NumberFormat = "@"
$c.Value = "2.59"
$c.Style = "Normal"
$ws.Range("E48").Value = "  +5.03%  "

$ws.Range("E49").Value = "  +2.14%  "

$c = $ws.Range("D50")
$c.NumberFormat = "@"
$c.Value = "133.69"
$c.Style = "Normal"
$ws.Range("E50").Value = "  +0.46%  "

$ws.Range("E51").Value = "  -0.04%  "
